# 16.6.1 worksheet update: replace forecast columns (2018 факт предв./2019 уточн./2020-2022 прогноз)
# with actual-data columns through 2021, adding a new утв./факт/% triplet for 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the formatting of row 2 (the thin divider row under the title) and
#        of the header row's last triplet (AJ3 style) into the three new columns
#        AK:AM before we populate them, so the new cells inherit the correct
#        borders/number formats instead of the plain column default.
$ws.Range("AJ2").Copy()
$ws.Range("AK2:AM2").PasteSpecial(-4122)

$ws.Range("AJ3").Copy()
$ws.Range("AK3:AM3").PasteSpecial(-4122)

$ws.Range("AJ12").Copy()
$ws.Range("AK12:AM12").PasteSpecial(-4122)

# --- 2. Header row 3: rename the forecast-era headers to actual-data headers and
#        add the new 2021 утв./факт/% triplet.
$ws.Range("AC3").Value = "2018 факт "
$ws.Range("AF3").Value = "2019 факт"
$ws.Range("AH3").Value = "2020 утв."
$ws.Range("AI3").Value = "2020 факт"
$ws.Range("AJ3").Value = "откл. от утв., %"
$ws.Range("AK3").Value = "2021 утв."
$ws.Range("AL3").Value = "2021 факт"
$ws.Range("AM3").Value = "откл. от утв., %"

# --- 3. Data rows: refresh the 2019 triplet (AE:AG), replace the old 2020/2021
#        forecast numbers in AH:AJ with 2020 actuals, and add the new 2021
#        утв./факт/% triplet in AK:AM.

# Row 5 - Государственные услуги общего назначения / General government services
$ws.Range("AE5").Value = 43737.8
$ws.Range("AF5").Value = 43258.3
$ws.Range("AG5").Value = 98.9
$ws.Range("AH5").Value = 46293.5
$ws.Range("AI5").Value = 47153.5
$ws.Range("AJ5").Value = 101.9
$ws.Range("AK5").Value = 47483.3
$ws.Range("AL5").Value = 52020.5
$ws.Range("AM5").Value = 109.6

# Row 6 - Экономические вопросы / Economic issues
$ws.Range("AE6").Value = 6265.4
$ws.Range("AF6").Value = 4434.6000000000004
$ws.Range("AG6").Value = 70.8
$ws.Range("AH6").Value = 7935.8
$ws.Range("AI6").Value = 3895.8
$ws.Range("AJ6").Value = 49.1
$ws.Range("AK6").Value = 8997
$ws.Range("AL6").Value = 6212.4
$ws.Range("AM6").Value = 69

# Row 7 - Охрана окружающей среды / Environmental protection
$ws.Range("AE7").Value = 728.5
$ws.Range("AF7").Value = 695.7
$ws.Range("AG7").Value = 95.5
$ws.Range("AH7").Value = 746.9
$ws.Range("AI7").Value = 583.20000000000005
$ws.Range("AJ7").Value = 78.099999999999994
$ws.Range("AK7").Value = 639.20000000000005
$ws.Range("AL7").Value = 600.79999999999995
$ws.Range("AM7").Value = 94

# Row 8 - Жилищные и коммунальные услуги / Housing and utilities services
$ws.Range("AE8").Value = 1249
$ws.Range("AF8").Value = 1244.7
$ws.Range("AG8").Value = 99.7
$ws.Range("AH8").Value = 1249
$ws.Range("AI8").Value = 1207.5999999999999
$ws.Range("AJ8").Value = 96.7
$ws.Range("AK8").Value = 1208.0999999999999
$ws.Range("AL8").Value = 1332.7
$ws.Range("AM8").Value = 110.3

# Row 9 - Здравоохранение / Healthcare
$ws.Range("AE9").Value = 2582.6
$ws.Range("AF9").Value = 2477.5
$ws.Range("AG9").Value = 95.9
$ws.Range("AH9").Value = 3109
$ws.Range("AI9").Value = 3225.2
$ws.Range("AJ9").Value = 103.7
$ws.Range("AK9").Value = 3131.3
$ws.Range("AL9").Value = 4833.7
$ws.Range("AM9").Value = 154.4

# Row 10 - Организация отдыха и культурно-религиозная деятельность / Recreation ...
$ws.Range("AE10").Value = 2686.4
$ws.Range("AF10").Value = 2829
$ws.Range("AG10").Value = 105.3
$ws.Range("AH10").Value = 2993.4
$ws.Range("AI10").Value = 2624.5
$ws.Range("AJ10").Value = 87.7
$ws.Range("AK10").Value = 2798.4
$ws.Range("AL10").Value = 3088
$ws.Range("AM10").Value = 110.3

# Row 11 - Образование / Education
$ws.Range("AE11").Value = 23397.4
$ws.Range("AF11").Value = 24364.799999999999
$ws.Range("AG11").Value = 104.1
$ws.Range("AH11").Value = 30085.9
$ws.Range("AI11").Value = 29223.5
$ws.Range("AJ11").Value = 97.1
$ws.Range("AK11").Value = 30439.7
$ws.Range("AL11").Value = 30705.3
$ws.Range("AM11").Value = 100.9

# Row 12 - Социальная защита / Social protection
$ws.Range("AE12").Value = 13137.1
$ws.Range("AF12").Value = 10924.7
$ws.Range("AG12").Value = 83.2
$ws.Range("AH12").Value = 12158.7
$ws.Range("AI12").Value = 10980.3
$ws.Range("AJ12").Value = 90.3
$ws.Range("AK12").Value = 11664.9
$ws.Range("AL12").Value = 11939.1
$ws.Range("AM12").Value = 102.4

# --- 4. Match the saved selection (cell AF4 is selected in the published file).
[void]$ws.Range("AF4").Select()
